# Updated cryptos list on Mon Feb 19 18:40:17 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Values are written as text (matching the
# original inline-string cells) so price strings like "52.080.43" and
# padded percentages like "  +0.72%  " are preserved verbatim.
$updates = [ordered]@{
    D2 = '52.080.43'
    E2 = '  +0.72%  '
    D3 = '2.938.74'
    E3 = '  +4.31%  '
    E4 = '  +0.00%  '
    D5 = '353.00'
    E5 = '  +0.49%  '
    D6 = '112.10'
    E6 = '  -1.03%  '
    D7 = '0.559'
    E7 = '  -0.24%  '
    E8 = '  +0.02%  '
    D9 = '0.625'
    E9 = '  +0.76%  '
    D10 = '39.42'
    E10 = '  -1.86%  '
    D11 = '0.0878'
    E11 = '  +3.40%  '
    E12 = '  +0.99%  '
    D13 = '20.10'
    E13 = '  +0.32%  '
    D14 = '3.403.87'
    E14 = '  +4.40%  '
    D15 = '7.77'
    E15 = '  -0.18%  '
    D16 = '2.938.68'
    E16 = '  +4.25%  '
    D17 = '0.981'
    E17 = '  +0.79%  '
    D18 = '52.127.39'
    E18 = '  +0.66%  '
    D19 = '7.62'
    E19 = '  +0.30%  '
    D20 = '3.29'
    E20 = '  -3.38%  '
    D21 = '14.21'
    E21 = '  +5.76%  '
    D22 = '0.0₃0979'
    E22 = '  +0.50%  '
    D23 = '71.18'
    E23 = '  +0.84%  '
    D24 = '268.45'
    E24 = '  -0.12%  '
    E25 = '  +0.56%  '
    E26 = '  +11.31%  '
    E27 = '  +2.97%  '
    E28 = '  -0.12%  '
    D29 = '7.13'
    E29 = '  +15.29%  '
    E30 = '  +15.14%  '
    D31 = '10.60'
    E31 = '  +0.50%  '
    E32 = '  -0.26%  '
    D33 = '37.03'
    E33 = '  -4.80%  '
    D34 = '6.10'
    E34 = '  +5.83%  '
    D35 = '53.04'
    E35 = '  +0.63%  '
    D36 = '0.0452'
    E36 = '  +0.25%  '
    E37 = '  -0.09%  '
    D38 = '3.38'
    E38 = '  +5.59%  '
    E39 = '  -2.92%  '
    E40 = '  +2.24%  '
    E41 = '  +4.31%  '
    E42 = '  +1.61%  '
    D43 = '23.47'
    E43 = '  +5.76%  '
    E44 = '  -1.61%  '
    E45 = '  +1.66%  '
    D46 = '2.199.06'
    E46 = '  +2.25%  '
    D47 = '3.52'
    E47 = '  +0.15%  '
    D48 = '112.24'
    E48 = '  -7.69%  '
    D49 = '0.249'
    E49 = '  +10.19%  '
    D50 = '0.0353'
    E50 = '  +9.43%  '
    D51 = '0.953'
    E51 = '  -3.93%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)

    # Price cells that look like a bare decimal number (e.g. "353.00",
    # "0.559") would otherwise be auto-coerced to a Number by COM,
    # silently dropping meaningful trailing zeros. Force those to stay
    # Text, same as the original inlineStr cell, using the classic
    # leading-apostrophe text-entry prefix, then restore the cell's
    # default (unstyled) appearance.
    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        $range.Value = "'" + $value
        $range.Style = 'Normal'
    } else {
        $range.Value = $value
    }
}
